# Apply edits to Tutorial 6 attendance sheet:
#  - Column A dates: change "/" separators to "-" separators (keep as text)
#  - For specific rows, flip D/E from 0->1 and H from 1->0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (format changed from dd/mm/yyyy to dd-mm-yyyy)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Rows where D (Total Attendance Count) and E (Real) go 0 -> 1,
# and H (Absent) goes 1 -> 0
$flipRows = @(4, 5, 11, 12, 13)

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)   # column A
    $cell.NumberFormat = "@"        # keep as text
    $cell.Value = $dates[$r]
}

foreach ($r in $flipRows) {
    $ws.Cells.Item($r, 4).Value = 1   # D - Total Attendance Count
    $ws.Cells.Item($r, 5).Value = 1   # E - Real
    $ws.Cells.Item($r, 8).Value = 0   # H - Absent
}
